# Apply updates to the "Intervention development" BCIO worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (agent) ---
$ws.Range("H2").Value = "independent continuant"
$ws.Range("T2").Value = 0

# --- Row 3 (behaviour change intervention development process) ---
$ws.Range("A3").Value = "BCIO:050277"
$ws.Range("D3").Value = "An intervention development process that is of a behaviour change intervention."

# --- Row 5 (collaboration process) ---
$ws.Range("D5").Value = "A development process in which agents work together to achieve some common objective."

# --- Row 6 (consultation -> consultation process) ---
$ws.Range("A6").Value = "BCIO:050278"
$ws.Range("B6").Value = "consultation process"
$ws.Range("D6").Value = "A development process in which a developer obtains relevant beliefs and experiences of stakeholder."

# --- Row 7 (developer) ---
$ws.Range("T7").Value = 0

# --- Row 8 (developer role) ---
$ws.Range("A8").Value = "BCIO:050279"
$ws.Range("D8").Value = "A role that involves having some responsibility for creating a product, service or intervention."

# --- Row 9 (development partner) ---
$ws.Range("D9").Value = "An agent that has a development partner role."
$ws.Range("H9").Value = "independent continuant"

# --- Row 10 (development partner role) ---
$ws.Range("A10").Value = "BCIO:050280"
$ws.Range("D10").Value = "A role that involves involves active participation in and a share of the responsibility for a development process."

# --- Row 12 (engagement -> intervention development process) ---
$ws.Range("A12").Value = "BCIO:050281"
$ws.Range("B12").Value = "intervention development process"
$ws.Range("D12").Value = "A development process that is for an intervention."

# --- Row 13 (intervention development process -> patient and public involvement) ---
$ws.Range("B13").Value = "patient and public involvement"

# --- Row 14 (patient and public involvement -> patient and public involvement and engagement) ---
$ws.Range("B14").Value = "patient and public involvement and engagement"

# --- Row 15 (patient and public involvement and engagement -> product development process) ---
$ws.Range("B15").Value = "product development process"

# --- Row 16 (product development process -> project development process) ---
$ws.Range("B16").Value = "project development process"

# --- Row 17 (project development process -> service development process) ---
$ws.Range("B17").Value = "service development process"

# --- Row 18 (service development process -> stakeholder) ---
$ws.Range("A18").Value = "BCIO:050276"
$ws.Range("B18").Value = "stakeholder"
$ws.Range("D18").Value = "An agent that has a stakeholder role."
$ws.Range("G18").Value = "agent"
$ws.Range("H18").Value = "independent continuant"
$ws.Range("M18").Value = "stakeholder role"

# --- Row 19 (stakeholder -> stakeholder engagement process) ---
$ws.Range("A19").Value = "BCIO:050282"
$ws.Range("B19").Value = "stakeholder engagement process"
$ws.Range("D19").Value = "A process in which in the course of development a developer discovers and takes account of the beliefs, feelings and experiences of stakeholders."
$ws.Range("G19").Value = "development process"
$ws.Range("H19").Value = "process"
$ws.Range("M19").Value = ""

# --- Row 20 (stakeholder role) ---
$ws.Range("H20").Value = "role"
